$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" '62.039.95'
$ws.Range("E2").Value = '  +2.84%  '

Set-TextValue $ws "D3" '3.413.48'
$ws.Range("E3").Value = '  +3.49%  '

Set-TextValue $ws "D4" '0.999'
$ws.Range("E4").Value = '  -0.05%  '

Set-TextValue $ws "D5" '577.67'
$ws.Range("E5").Value = '  +2.92%  '

Set-TextValue $ws "D6" '138.83'
$ws.Range("E6").Value = '  +7.51%  '

Set-TextValue $ws "D8" '3.412.91'
$ws.Range("E8").Value = '  +3.45%  '

$ws.Range("E9").Value = '  +1.34%  '

Set-TextValue $ws "D10" '7.53'
$ws.Range("E10").Value = '  +2.39%  '

$ws.Range("E11").Value = '  +9.73%  '

Set-TextValue $ws "D12" '0.397'
$ws.Range("E12").Value = '  +6.98%  '

Set-TextValue $ws "D13" '3.994.75'
$ws.Range("E13").Value = '  +3.44%  '

$ws.Range("E14").Value = '  +1.95%  '

$ws.Range("E15").Value = '  +8.56%  '

Set-TextValue $ws "D16" '3.415.95'
$ws.Range("E16").Value = '  +3.53%  '

Set-TextValue $ws "D17" '25.56'
$ws.Range("E17").Value = '  +5.78%  '

Set-TextValue $ws "D18" '61.993.26'
$ws.Range("E18").Value = '  +2.48%  '

Set-TextValue $ws "D19" '14.17'
$ws.Range("E19").Value = '  +6.45%  '

Set-TextValue $ws "D20" '5.93'
$ws.Range("E20").Value = '  +4.77%  '

Set-TextValue $ws "D21" '9.50'
$ws.Range("E21").Value = '  +6.16%  '

Set-TextValue $ws "D22" '391.57'
$ws.Range("E22").Value = '  +11.73%  '

Set-TextValue $ws "D23" '0.575'
$ws.Range("E23").Value = '  +3.89%  '

Set-TextValue $ws "D24" '3.549.21'
$ws.Range("E24").Value = '  +3.54%  '

$ws.Range("E25").Value = '  +19.69%  '

$ws.Range("E26").Value = '  +0.23%  '

Set-TextValue $ws "D27" '71.65'
$ws.Range("E27").Value = '  +3.57%  '

Set-TextValue $ws "D28" '1.59'
$ws.Range("E28").Value = '  +10.50%  '

$ws.Range("E29").Value = '  +4.73%  '

Set-TextValue $ws "D30" '0.998'
$ws.Range("E30").Value = '  -0.21%  '

Set-TextValue $ws "D31" '8.34'
$ws.Range("E31").Value = '  +6.71%  '

Set-TextValue $ws "D32" '0.160'
$ws.Range("E32").Value = '  +5.51%  '

$ws.Range("E33").Value = '  +3.54%  '

Set-TextValue $ws "D34" '3.442.46'
$ws.Range("E34").Value = '  +3.55%  '

$ws.Range("E35").Value = '  -0.03%  '

Set-TextValue $ws "D36" '23.64'
$ws.Range("E36").Value = '  +4.51%  '

Set-TextValue $ws "D37" '5.56'
$ws.Range("E37").Value = '  +5.78%  '

Set-TextValue $ws "D38" '7.01'
$ws.Range("E38").Value = '  +3.75%  '

Set-TextValue $ws "D39" '1.56'
$ws.Range("E39").Value = '  +5.81%  '

Set-TextValue $ws "D40" '161.59'
$ws.Range("E40").Value = '  +2.58%  '

Set-TextValue $ws "D41" '0.0797'
$ws.Range("E41").Value = '  +6.19%  '

Set-TextValue $ws "D42" '1.76'
$ws.Range("E42").Value = '  +14.38%  '

$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("E44").Value = '  +6.71%  '

Set-TextValue $ws "D45" '0.779'
$ws.Range("E45").Value = '  +5.08%  '

Set-TextValue $ws "D46" '4.49'
$ws.Range("E46").Value = '  +3.70%  '

Set-TextValue $ws "D47" '25.20'
$ws.Range("E47").Value = '  +10.85%  '

$ws.Range("E48").Value = '  +1.72%  '

Set-TextValue $ws "D49" '7.00'
$ws.Range("E49").Value = '  +4.99%  '

Set-TextValue $ws "D50" '23.00'
$ws.Range("E50").Value = '  +6.60%  '

Set-TextValue $ws "D51" '2.397.08'
$ws.Range("E51").Value = '  +10.85%  '
